$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expected concatenation results per row (G column), matching the
# "10"+20-as-various-types test matrix. Values are set as TEXT (quote
# prefixed, Text number format) rather than numbers.
$results = @{
    5  = "1020"
    6  = "1020"
    7  = "1020"
    8  = "1020"
    9  = "1020.0"
    10 = "1020.0"
    11 = "1020"
    12 = "1020"
    13 = "1020"
    14 = "1020"
    15 = "10.020"
    16 = "10.020"
    17 = "1020"
    18 = "1020"
    19 = "1020"
    20 = "1020"
    21 = "1020.0"
    22 = "1020.0"
    23 = "1020"
    24 = "1020"
    25 = "1020"
    26 = "1020"
    27 = "10.020"
    28 = "10.020"
    29 = "1020"
    30 = "1020"
    31 = "1020"
    32 = "1020"
    33 = "1020.0"
    34 = "1020.0"
    35 = "1020"
    36 = "1020"
    37 = "1020"
    38 = "1020"
    39 = "10.020"
    40 = "10.020"
}

foreach ($r in 5..40) {
    $cell = $ws.Range("G" + $r)
    $cell.Value = "'" + $results[$r]
    $cell.NumberFormat = "@"
}

# G41 already held the text "1020" - leave its value, just make sure it
# keeps the quoted-text formatting consistent with the rest of the column.
$g41 = $ws.Range("G41")
$g41.Value = "'1020"
$g41.NumberFormat = "@"

# Last selection made while editing, stored in the saved sheet view.
$ws.Range("L14").Select()
